# Insert a new data row at row 171 (pushes existing rows 171.. down by one)
# then populate it with the required values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 171 so everything from old row 171
# onward shifts down by one (old row 171 becomes row 172, ..., old row 232
# becomes row 233).
$ws.Rows.Item(171).Insert()

# Populate the newly-inserted row 171 with the same "shape" of data as the
# row that used to occupy this slot (now at row 172), changing only the
# price/volume fields that differ for this particular record.
$ws.Cells.Item(171, 1).Value = 11
$ws.Cells.Item(171, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(171, 3).Value = "Bíobío"
$ws.Cells.Item(171, 4).Value = 45146
$ws.Cells.Item(171, 5).Value = 8
$ws.Cells.Item(171, 6).Value = 100112043
$ws.Cells.Item(171, 7).Value = "Pepino ensalada"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 100
$ws.Cells.Item(171, 11).Value = 11000
$ws.Cells.Item(171, 12).Value = 12000
$ws.Cells.Item(171, 13).Value = 11500
$ws.Cells.Item(171, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(171, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(171, 16).Value = 192
$ws.Cells.Item(171, 17).Value = 60
$ws.Cells.Item(171, 18).Value = "Hortaliza"
